$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value. All values must remain TEXT
# (matching the source inlineStr cells), so NumberFormat is forced to
# "@" (Text) before assignment to stop Excel from auto-coercing
# numeric-looking strings (e.g. "255.32") or percentages (e.g. "3.79%")
# into numbers.
$updates = @{
    'D2' = '255.32'
    'E2' = '3.79%'
    'D3' = '28.16'
    'E3' = '-5.51%'
    'D4' = '5.260'
    'E4' = '2.09%'
    'E5' = '1.47%'
    'D6' = '6.706'
    'E6' = '0.78%'
    'D7' = '0.8680'
    'E7' = '2.10%'
    'D8' = '1.049'
    'E8' = '22.91%'
    'D9' = '0.1414'
    'E9' = '2.51%'
    'D10' = '0.07117'
    'E10' = '0.51%'
    'D11' = '0.03173'
    'E11' = '-2.61%'
    'D12' = '0.09226'
    'E12' = '-1.51%'
    'D13' = '0.001539'
    'E13' = '0.14%'
    'D14' = '0.0006091'
    'D15' = '0.005806'
    'E15' = '-3.31%'
    'D16' = '3.498'
    'E16' = '-0.46%'
    'E17' = '-0.35%'
    'E18' = '-0.89%'
    'D19' = '0.3178'
    'E19' = '0.45%'
    'D20' = '0.03468'
    'E20' = '2.72%'
    'E21' = '0.04%'
    'D22' = '3.562'
    'E22' = '2.31%'
    'D23' = '0.04145'
    'E23' = '0.26%'
    'E24' = '-4.56%'
    'D25' = '0.001226'
    'E25' = '-0.12%'
    'D26' = '0.004880'
    'E26' = '17.82%'
    'E27' = '0.03%'
    'D28' = '0.00008001'
    'E28' = '-44.76%'
    'D40' = '0.03820'
    'E40' = '2.03%'
    'B41' = 'KickToken'
    'C41' = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
    'D41' = '0.005761'
    'E41' = '0.21%'
    'B42' = 'BKEXToken'
    'C42' = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
    'D42' = '0.1102'
    'E42' = '2.90%'
    'D43' = '0.002343'
    'E43' = '1.91%'
    'D44' = '0.01079'
    'E44' = '25.81%'
    'E45' = '-4.97%'
    'E46' = '0.03%'
    'D47' = '0.09302'
    'E47' = '31.05%'
    'D48' = '0.002152'
    'E48' = '-3.57%'
    'E49' = '0.03%'
    'E50' = '0.03%'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
